# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# Rebuilds the 16 worker/period data rows (rows 16-31) of the "Hoja1"
# worksheet: both workers now carry all 8 periods (2205..2110, descending),
# and "Salario Basico" (col G) for LUIS MIGUEL MEZA OJEDA drops from
# 1,000,000 to 908,526.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# CC / 73118385 / LUIS MIGUEL MEZA OJEDA block (rows 16-23)
$rowsMeza = @(
    @{ Row = 16; Periodo = "2205"; ValorMora = 27861; Salario = 908526 },
    @{ Row = 17; Periodo = "2204"; ValorMora = 36341; Salario = 908526 },
    @{ Row = 18; Periodo = "2203"; ValorMora = 36341; Salario = 908526 },
    @{ Row = 19; Periodo = "2202"; ValorMora = 36341; Salario = 908526 },
    @{ Row = 20; Periodo = "2201"; ValorMora = 36341; Salario = 908526 },
    @{ Row = 21; Periodo = "2112"; ValorMora = 36341; Salario = 908526 },
    @{ Row = 22; Periodo = "2111"; ValorMora = 36341; Salario = 908526 },
    @{ Row = 23; Periodo = "2110"; ValorMora = 36341; Salario = 908526 }
)

foreach ($item in $rowsMeza) {
    $r = $item.Row
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = "73118385"
    $ws.Range("D$r").Value = "LUIS MIGUEL MEZA OJEDA"
    $ws.Range("E$r").Value = $item.Periodo
    $ws.Range("F$r").Value = $item.ValorMora
    $ws.Range("G$r").Value = $item.Salario
}

# CE / 488479 / JOSE LUENGO NAVARRO block (rows 24-31)
$rowsLuengo = @(
    @{ Row = 24; Periodo = "2205"; ValorMora = 46000; Salario = 1500000 },
    @{ Row = 25; Periodo = "2204"; ValorMora = 60000; Salario = 1500000 },
    @{ Row = 26; Periodo = "2203"; ValorMora = 60000; Salario = 1500000 },
    @{ Row = 27; Periodo = "2202"; ValorMora = 60000; Salario = 1500000 },
    @{ Row = 28; Periodo = "2201"; ValorMora = 60000; Salario = 1500000 },
    @{ Row = 29; Periodo = "2112"; ValorMora = 60000; Salario = 1500000 },
    @{ Row = 30; Periodo = "2111"; ValorMora = 60000; Salario = 1500000 },
    @{ Row = 31; Periodo = "2110"; ValorMora = 60000; Salario = 1500000 }
)

foreach ($item in $rowsLuengo) {
    $r = $item.Row
    $ws.Range("B$r").Value = "CE"
    $ws.Range("C$r").Value = "488479"
    $ws.Range("D$r").Value = "JOSE LUENGO NAVARRO"
    $ws.Range("E$r").Value = $item.Periodo
    $ws.Range("F$r").Value = $item.ValorMora
    $ws.Range("G$r").Value = $item.Salario
}
